$d = $word.ActiveDocument

# Update the date heading.
$d.Content.Find.Execute("2023-05-12 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-13 Saturday", 2) | Out-Null

# Update each arithmetic-problem cell in the table. Every "old" string below
# is unique within the document, and none of the "new" strings collide with
# any "old" string, so the replacements are order-independent and safe to
# run as a straight sequence of literal Find/Replace operations.
$d.Content.Find.Execute("98-15=", $true, $false, $false, $false, $false, $true, 1, $false, "1+36=", 2) | Out-Null
$d.Content.Find.Execute("30+26=", $true, $false, $false, $false, $false, $true, 1, $false, "79-0=", 2) | Out-Null
$d.Content.Find.Execute("43-31=", $true, $false, $false, $false, $false, $true, 1, $false, "29-21=", 2) | Out-Null
$d.Content.Find.Execute("8+31=", $true, $false, $false, $false, $false, $true, 1, $false, "81-48=", 2) | Out-Null
$d.Content.Find.Execute("6+38=", $true, $false, $false, $false, $false, $true, 1, $false, "34+26=", 2) | Out-Null
$d.Content.Find.Execute("9+11=", $true, $false, $false, $false, $false, $true, 1, $false, "53-29=", 2) | Out-Null
$d.Content.Find.Execute("76-61=", $true, $false, $false, $false, $false, $true, 1, $false, "37+61=", 2) | Out-Null
$d.Content.Find.Execute("69-49=", $true, $false, $false, $false, $false, $true, 1, $false, "96-54=", 2) | Out-Null
$d.Content.Find.Execute("34+57=", $true, $false, $false, $false, $false, $true, 1, $false, "7+46=", 2) | Out-Null
$d.Content.Find.Execute("35+21=", $true, $false, $false, $false, $false, $true, 1, $false, "51-26=", 2) | Out-Null
$d.Content.Find.Execute("35-24=", $true, $false, $false, $false, $false, $true, 1, $false, "85-57=", 2) | Out-Null
$d.Content.Find.Execute("68-13=", $true, $false, $false, $false, $false, $true, 1, $false, "5+35=", 2) | Out-Null
$d.Content.Find.Execute("1+78=", $true, $false, $false, $false, $false, $true, 1, $false, "24+29=", 2) | Out-Null
$d.Content.Find.Execute("69-27=", $true, $false, $false, $false, $false, $true, 1, $false, "42-29=", 2) | Out-Null
$d.Content.Find.Execute("36-7=", $true, $false, $false, $false, $false, $true, 1, $false, "23+40=", 2) | Out-Null
$d.Content.Find.Execute("39+29=", $true, $false, $false, $false, $false, $true, 1, $false, "87-36=", 2) | Out-Null
$d.Content.Find.Execute("88-41=", $true, $false, $false, $false, $false, $true, 1, $false, "23+22=", 2) | Out-Null
$d.Content.Find.Execute("73+11=", $true, $false, $false, $false, $false, $true, 1, $false, "29+41=", 2) | Out-Null
$d.Content.Find.Execute("88-51=", $true, $false, $false, $false, $false, $true, 1, $false, "82-64=", 2) | Out-Null
$d.Content.Find.Execute("40-20=", $true, $false, $false, $false, $false, $true, 1, $false, "32-5=", 2) | Out-Null
$d.Content.Find.Execute("89-76=", $true, $false, $false, $false, $false, $true, 1, $false, "69+28=", 2) | Out-Null
$d.Content.Find.Execute("35+7=", $true, $false, $false, $false, $false, $true, 1, $false, "13-3=", 2) | Out-Null
$d.Content.Find.Execute("92-19=", $true, $false, $false, $false, $false, $true, 1, $false, "59+39=", 2) | Out-Null
$d.Content.Find.Execute("11+43=", $true, $false, $false, $false, $false, $true, 1, $false, "56+31=", 2) | Out-Null
$d.Content.Find.Execute("5+71=", $true, $false, $false, $false, $false, $true, 1, $false, "76-33=", 2) | Out-Null
$d.Content.Find.Execute("96-25=", $true, $false, $false, $false, $false, $true, 1, $false, "95-69=", 2) | Out-Null
$d.Content.Find.Execute("36+37=", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=", 2) | Out-Null
$d.Content.Find.Execute("96-39=", $true, $false, $false, $false, $false, $true, 1, $false, "36-8=", 2) | Out-Null
$d.Content.Find.Execute("26-7=", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=", 2) | Out-Null
$d.Content.Find.Execute("0+97=", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=", 2) | Out-Null
$d.Content.Find.Execute("99-23=", $true, $false, $false, $false, $false, $true, 1, $false, "81-23=", 2) | Out-Null
$d.Content.Find.Execute("50+18=", $true, $false, $false, $false, $false, $true, 1, $false, "31+4=", 2) | Out-Null
$d.Content.Find.Execute("26+35=", $true, $false, $false, $false, $false, $true, 1, $false, "56-16=", 2) | Out-Null
$d.Content.Find.Execute("60+27=", $true, $false, $false, $false, $false, $true, 1, $false, "64+3=", 2) | Out-Null
$d.Content.Find.Execute("60-34=", $true, $false, $false, $false, $false, $true, 1, $false, "83-21=", 2) | Out-Null
$d.Content.Find.Execute("35-27=", $true, $false, $false, $false, $false, $true, 1, $false, "2+22=", 2) | Out-Null
$d.Content.Find.Execute("11+23=", $true, $false, $false, $false, $false, $true, 1, $false, "35+22=", 2) | Out-Null
$d.Content.Find.Execute("31-15=", $true, $false, $false, $false, $false, $true, 1, $false, "31-20=", 2) | Out-Null
$d.Content.Find.Execute("94-68=", $true, $false, $false, $false, $false, $true, 1, $false, "78-76=", 2) | Out-Null
$d.Content.Find.Execute("81-19=", $true, $false, $false, $false, $false, $true, 1, $false, "86+10=", 2) | Out-Null
$d.Content.Find.Execute("57-51=", $true, $false, $false, $false, $false, $true, 1, $false, "96-91=", 2) | Out-Null
$d.Content.Find.Execute("95-0=", $true, $false, $false, $false, $false, $true, 1, $false, "67-61=", 2) | Out-Null
$d.Content.Find.Execute("98+0=", $true, $false, $false, $false, $false, $true, 1, $false, "43+22=", 2) | Out-Null
$d.Content.Find.Execute("33+58=", $true, $false, $false, $false, $false, $true, 1, $false, "40-23=", 2) | Out-Null
$d.Content.Find.Execute("41-22=", $true, $false, $false, $false, $false, $true, 1, $false, "21+45=", 2) | Out-Null
$d.Content.Find.Execute("99-74=", $true, $false, $false, $false, $false, $true, 1, $false, "18+71=", 2) | Out-Null
$d.Content.Find.Execute("36+4=", $true, $false, $false, $false, $false, $true, 1, $false, "58-0=", 2) | Out-Null
$d.Content.Find.Execute("6+23=", $true, $false, $false, $false, $false, $true, 1, $false, "66+16=", 2) | Out-Null
$d.Content.Find.Execute("81-79=", $true, $false, $false, $false, $false, $true, 1, $false, "10+3=", 2) | Out-Null
$d.Content.Find.Execute("47-21=", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=", 2) | Out-Null
$d.Content.Find.Execute("67-56=", $true, $false, $false, $false, $false, $true, 1, $false, "42-35=", 2) | Out-Null
$d.Content.Find.Execute("40+13=", $true, $false, $false, $false, $false, $true, 1, $false, "57-27=", 2) | Out-Null
$d.Content.Find.Execute("68-39=", $true, $false, $false, $false, $false, $true, 1, $false, "68-0=", 2) | Out-Null
$d.Content.Find.Execute("57-56=", $true, $false, $false, $false, $false, $true, 1, $false, "76-25=", 2) | Out-Null
$d.Content.Find.Execute("84-55=", $true, $false, $false, $false, $false, $true, 1, $false, "68-66=", 2) | Out-Null
$d.Content.Find.Execute("6+21=", $true, $false, $false, $false, $false, $true, 1, $false, "89-20=", 2) | Out-Null
$d.Content.Find.Execute("85-72=", $true, $false, $false, $false, $false, $true, 1, $false, "29+6=", 2) | Out-Null
$d.Content.Find.Execute("18-17=", $true, $false, $false, $false, $false, $true, 1, $false, "53-49=", 2) | Out-Null
$d.Content.Find.Execute("21+54=", $true, $false, $false, $false, $false, $true, 1, $false, "17+76=", 2) | Out-Null
$d.Content.Find.Execute("38+20=", $true, $false, $false, $false, $false, $true, 1, $false, "59+1=", 2) | Out-Null
$d.Content.Find.Execute("52-27=", $true, $false, $false, $false, $false, $true, 1, $false, "28+30=", 2) | Out-Null
$d.Content.Find.Execute("44-7=", $true, $false, $false, $false, $false, $true, 1, $false, "18-2=", 2) | Out-Null
$d.Content.Find.Execute("18+21=", $true, $false, $false, $false, $false, $true, 1, $false, "0+3=", 2) | Out-Null
$d.Content.Find.Execute("15+38=", $true, $false, $false, $false, $false, $true, 1, $false, "16+17=", 2) | Out-Null
$d.Content.Find.Execute("22-18=", $true, $false, $false, $false, $false, $true, 1, $false, "28-12=", 2) | Out-Null
$d.Content.Find.Execute("3+53=", $true, $false, $false, $false, $false, $true, 1, $false, "42+40=", 2) | Out-Null
$d.Content.Find.Execute("94-62=", $true, $false, $false, $false, $false, $true, 1, $false, "5+14=", 2) | Out-Null
$d.Content.Find.Execute("30+48=", $true, $false, $false, $false, $false, $true, 1, $false, "55-45=", 2) | Out-Null
$d.Content.Find.Execute("40-36=", $true, $false, $false, $false, $false, $true, 1, $false, "11+40=", 2) | Out-Null
$d.Content.Find.Execute("9+72=", $true, $false, $false, $false, $false, $true, 1, $false, "44+48=", 2) | Out-Null
$d.Content.Find.Execute("93-86=", $true, $false, $false, $false, $false, $true, 1, $false, "69-29=", 2) | Out-Null
$d.Content.Find.Execute("38-16=", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=", 2) | Out-Null
$d.Content.Find.Execute("36-3=", $true, $false, $false, $false, $false, $true, 1, $false, "71+7=", 2) | Out-Null
$d.Content.Find.Execute("39-2=", $true, $false, $false, $false, $false, $true, 1, $false, "27+41=", 2) | Out-Null
$d.Content.Find.Execute("88-35=", $true, $false, $false, $false, $false, $true, 1, $false, "98-76=", 2) | Out-Null
$d.Content.Find.Execute("40-29=", $true, $false, $false, $false, $false, $true, 1, $false, "78-47=", 2) | Out-Null
$d.Content.Find.Execute("86+4=", $true, $false, $false, $false, $false, $true, 1, $false, "38+19=", 2) | Out-Null
$d.Content.Find.Execute("61-32=", $true, $false, $false, $false, $false, $true, 1, $false, "17+63=", 2) | Out-Null
$d.Content.Find.Execute("21+77=", $true, $false, $false, $false, $false, $true, 1, $false, "75+23=", 2) | Out-Null
$d.Content.Find.Execute("84-82=", $true, $false, $false, $false, $false, $true, 1, $false, "43-9=", 2) | Out-Null
$d.Content.Find.Execute("27+38=", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=", 2) | Out-Null
$d.Content.Find.Execute("56-18=", $true, $false, $false, $false, $false, $true, 1, $false, "50+28=", 2) | Out-Null
$d.Content.Find.Execute("55+14=", $true, $false, $false, $false, $false, $true, 1, $false, "67-62=", 2) | Out-Null
$d.Content.Find.Execute("59-35=", $true, $false, $false, $false, $false, $true, 1, $false, "82-76=", 2) | Out-Null
$d.Content.Find.Execute("97-5=", $true, $false, $false, $false, $false, $true, 1, $false, "6+69=", 2) | Out-Null
$d.Content.Find.Execute("9+51=", $true, $false, $false, $false, $false, $true, 1, $false, "45-17=", 2) | Out-Null
$d.Content.Find.Execute("52-11=", $true, $false, $false, $false, $false, $true, 1, $false, "85-52=", 2) | Out-Null
$d.Content.Find.Execute("58-58=", $true, $false, $false, $false, $false, $true, 1, $false, "74+4=", 2) | Out-Null
$d.Content.Find.Execute("2+67=", $true, $false, $false, $false, $false, $true, 1, $false, "36-27=", 2) | Out-Null
$d.Content.Find.Execute("38-6=", $true, $false, $false, $false, $false, $true, 1, $false, "35-20=", 2) | Out-Null
$d.Content.Find.Execute("11+46=", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=", 2) | Out-Null
$d.Content.Find.Execute("92-43=", $true, $false, $false, $false, $false, $true, 1, $false, "7+75=", 2) | Out-Null
$d.Content.Find.Execute("9+68=", $true, $false, $false, $false, $false, $true, 1, $false, "38+49=", 2) | Out-Null
$d.Content.Find.Execute("18+76=", $true, $false, $false, $false, $false, $true, 1, $false, "89-63=", 2) | Out-Null
$d.Content.Find.Execute("28+27=", $true, $false, $false, $false, $false, $true, 1, $false, "79-51=", 2) | Out-Null
$d.Content.Find.Execute("18+56=", $true, $false, $false, $false, $false, $true, 1, $false, "14+10=", 2) | Out-Null
$d.Content.Find.Execute("84+6=", $true, $false, $false, $false, $false, $true, 1, $false, "80-73=", 2) | Out-Null
$d.Content.Find.Execute("6+43=", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=", 2) | Out-Null
$d.Content.Find.Execute("36+50=", $true, $false, $false, $false, $false, $true, 1, $false, "54+26=", 2) | Out-Null
$d.Content.Find.Execute("96-29=", $true, $false, $false, $false, $false, $true, 1, $false, "46-11=", 2) | Out-Null
